# Leave card update: add a new VL(2-0-0) leave record for 12/27,28/2022
# under the 2022 section, add a "2023" year header row, backfill the
# month-start dates for the (previously blank) PERIOD column for the
# following ~4 years of rows, and record two EARNED entries (1.25 each)
# plus recompute the resulting balances. Also refresh the footer
# signatory block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1) Grow the Table1 list object by one row (table range A8:K139 -> A8:K140)
# ---------------------------------------------------------------------
$earnedFormula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# Preserve the special "last row" formatting (rows 139 -> 140) before we
# turn the old last row back into a normal data row.
$ws.Range("A139:K139").Copy($ws.Range("A140:K140"))
$ws.Range("G140").Formula = $earnedFormula

# Old row 139 becomes an ordinary data row (copy format from row 138).
$ws.Range("A138:K138").Copy($ws.Range("A139:K139"))
$ws.Range("G139").Formula = $earnedFormula

# Extend the ListObject / table definition to match the new extent.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A8:K140"))

# ---------------------------------------------------------------------
# 2) New year header row (87): "2023" (text, matches other year rows)
#    Written first so it lands on the same shared-string index ordering
#    as the authoring app used.
# ---------------------------------------------------------------------
$ws.Range("A87").NumberFormat = "@"
$ws.Range("A87").Value = "2023"
$ws.Range("A10").Copy()
$ws.Range("A87").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3) New leave entry row (86): VL(2-0-0), 2 days, taken 12/27,28/2022
# ---------------------------------------------------------------------
$ws.Range("B86").Value = "VL(2-0-0)"
$ws.Range("D86").Value = 2

$ws.Range("K30").Copy()
$ws.Range("K86").PasteSpecial(-4122)
$ws.Range("K86").Value = "12/27,28/2022"

# ---------------------------------------------------------------------
# 4) Backfill PERIOD (month-start) dates for rows 88-134
# ---------------------------------------------------------------------
$periodDates = @{
  88  = 44927; 89  = 44958; 90  = 44986; 91  = 45017; 92  = 45047;
  93  = 45078; 94  = 45108; 95  = 45139; 96  = 45170; 97  = 45200;
  98  = 45231; 99  = 45261; 100 = 45292; 101 = 45323; 102 = 45352;
  103 = 45383; 104 = 45413; 105 = 45444; 106 = 45474; 107 = 45505;
  108 = 45536; 109 = 45566; 110 = 45597; 111 = 45627; 112 = 45658;
  113 = 45689; 114 = 45717; 115 = 45748; 116 = 45778; 117 = 45809;
  118 = 45839; 119 = 45870; 120 = 45901; 121 = 45931; 122 = 45962;
  123 = 45992; 124 = 46023; 125 = 46054; 126 = 46082; 127 = 46113;
  128 = 46143; 129 = 46174; 130 = 46204; 131 = 46235; 132 = 46266;
  133 = 46296; 134 = 46327
}
foreach ($row in $periodDates.Keys) {
  $ws.Range("A$row").Value = $periodDates[$row]
}

# ---------------------------------------------------------------------
# 5) New EARNED entries: 1.25 on rows 88 and 89 (EARNED column C)
# ---------------------------------------------------------------------
$ws.Range("C88").Value = 1.25
$ws.Range("C89").Value = 1.25

# ---------------------------------------------------------------------
# 6) View bookkeeping: keep the split pane + selections pointed at the
#    newly edited rows.
# ---------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.SplitRow = 76
$ws.Range("K87").Select()

# ---------------------------------------------------------------------
# 7) Footer signatory update
# ---------------------------------------------------------------------
$ps = $ws.PageSetup
$ps.CenterFooter = "`nCERTIFIED CORRECT BY: &UNANETTE B. SUSA&U`n                                              OIC - HRMO"

$wb.Application.Calculate()
